$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.030.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.430.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.797.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.431.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.029.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "

$ws.Range("E31").Value = "  +18.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0758"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "130.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +23.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.950.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.56%  "

$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.659.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
